# Generate Report for Handback
# Updates the timestamp strings recorded on the "Overview", "zh-cn" and
# "de-de" sheets to reflect a newer handoff/handback generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# G2 = "Latest HO Xliff Generate Date" for the first row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-15 18:59:10"

# --- zh-cn sheet ------------------------------------------------------
# H2 = "Correspond Handoff Datetime", K2 = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-15 18:58:58"
$wsZhCn.Range("K2").Value = "2016-08-15 18:59:29"

# --- de-de sheet --------------------------------------------------------
# H2 = "Correspond Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-15 18:59:36"
